$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price & volume changes)
# Force text format only where the new value would otherwise be
# auto-converted to a number, so formatting (e.g. trailing zeros) is kept

$ws.Range('D2').Value = '26.199.63'
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('D3').Value = '1.659.75'
$ws.Range('E3').Value = '  -1.49%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.24'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5219'
$ws.Range('E6').Value = '  -2.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2665'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06317'
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('E10').Value = '  -3.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07716'
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('D12').Value = '1.662.61'
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('E13').Value = '  -1.91%  '
$ws.Range('D14').Value = '1.887.49'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('E15').Value = '  -2.69%  '
$ws.Range('D16').Value = '0.0₅8206'
$ws.Range('E16').Value = '  -2.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.77'
$ws.Range('E17').Value = '  -1.88%  '
$ws.Range('D18').Value = '26.231.28'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('E20').Value = '  -2.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.45'
$ws.Range('E21').Value = '  -1.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.14'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.078'
$ws.Range('E23').Value = '  -4.58%  '
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '138.39'
$ws.Range('E25').Value = '  -4.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1237'
$ws.Range('E26').Value = '  -3.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.220'
$ws.Range('E27').Value = '  -3.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.13'
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('E29').Value = '  -1.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05998'
$ws.Range('E30').Value = '  -2.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.280'
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.626'
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.310'
$ws.Range('E33').Value = '  -4.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.635'
$ws.Range('E34').Value = '  -3.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9778'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.780'
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5889'
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01586'
$ws.Range('E39').Value = '  -3.64%  '
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8646'
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').Value = '1.031.52'
$ws.Range('E43').Value = '  -3.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.55'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('D45').Value = '1.801.87'
$ws.Range('E45').Value = '  -1.87%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.00'
$ws.Range('E47').Value = '  -0.51%  '
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.084'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05185'
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4230'
$ws.Range('E51').Value = '  -0.31%  '

Write-Host "Updated cryptos list"
